$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simplify the steel description text in B2: remove the "/RME" segment
$cell = $ws.Range("B2")
$cell.Value = "25% S/LFBR+CDH/H:2`n20% S/LFM+CDM/H:1`n5% S+SL/LFM+CDL/H:1`n38% CR/LFM+CDM/H:1`n12% CR+PC/LFM+CDL/H:2"

# Wrap the text and grow the row so the full description is visible
$cell.WrapText = $true
$ws.Rows.Item(2).RowHeight = 256

# Leave the selection where the author last left it when saving
$null = $ws.Range("B9").Select()
